# Scheduled runner refresh: update market-price-derived columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 241.33333
$ws.Range("I12").Value = 249.33333
$ws.Range("J12").Value = 233.33333
$ws.Range("K12").Value = 249.33333
$ws.Range("L12").Value = 233.33333
$ws.Range("M12").Value = -79.33332999999999
$ws.Range("N12").Value = -573.3333299999999
$ws.Range("H28").Value = 532.5
$ws.Range("I28").Value = 532.5
$ws.Range("K28").Value = 532.5
$ws.Range("M28").Value = -47.5
$ws.Range("H98").Value = 1003
$ws.Range("I98").Value = 917.4286
$ws.Range("J98").Value = 1174.1428
$ws.Range("K98").Value = 917.4286
$ws.Range("L98").Value = 1174.1428
$ws.Range("M98").Value = 580.5714
$ws.Range("N98").Value = -4170.1428
$ws.Range("H112").Value = 2168.3333
$ws.Range("J112").Value = 2562
$ws.Range("L112").Value = 7686
$ws.Range("N112").Value = -9902
$ws.Range("H122").Value = 1003
$ws.Range("I122").Value = 917.4286
$ws.Range("J122").Value = 1174.1428
$ws.Range("K122").Value = 2752.2858
$ws.Range("L122").Value = 3522.4284
$ws.Range("M122").Value = -302.2857999999997
$ws.Range("N122").Value = -8422.428400000001
$ws.Range("H129").Value = 983.62195
$ws.Range("I129").Value = 333.33334
$ws.Range("J129").Value = 1034.9606
$ws.Range("K129").Value = 1000.00002
$ws.Range("L129").Value = 3104.8818
$ws.Range("M129").Value = 3999.99998
$ws.Range("N129").Value = -13104.8818
$ws.Range("H132").Value = 2254.9583
$ws.Range("I132").Value = 2242.5264
$ws.Range("J132").Value = 2302.2
$ws.Range("K132").Value = 6727.5792
$ws.Range("L132").Value = 6906.599999999999
$ws.Range("M132").Value = -4197.5792
$ws.Range("N132").Value = -11966.6
$ws.Range("H135").Value = 782
$ws.Range("I135").Value = 728
$ws.Range("J135").Value = 998
$ws.Range("K135").Value = 6552
$ws.Range("L135").Value = 8982
$ws.Range("M135").Value = -4017
$ws.Range("N135").Value = -14052
$ws.Range("H138").Value = 3413.8955
$ws.Range("I138").Value = 755.75
$ws.Range("J138").Value = 4897.5117
$ws.Range("K138").Value = 2267.25
$ws.Range("L138").Value = 14692.5351
$ws.Range("M138").Value = 2872.75
$ws.Range("N138").Value = -24972.5351

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3457.68
$ws.Range("I45").Value = 3305.2354
$ws.Range("J45").Value = 3781.625
$ws.Range("K45").Value = 3305.2354
$ws.Range("L45").Value = 3781.625
$ws.Range("M45").Value = -2928.2354
$ws.Range("N45").Value = -4535.625
$ws.Range("H61").Value = 2476.9524
$ws.Range("I61").Value = 1393
$ws.Range("J61").Value = 4644.857
$ws.Range("K61").Value = 1393
$ws.Range("L61").Value = 4644.857
$ws.Range("M61").Value = -1181
$ws.Range("N61").Value = -5068.857
$ws.Range("H74").Value = 4049.9697
$ws.Range("I74").Value = 4269.6665
$ws.Range("J74").Value = 1853
$ws.Range("K74").Value = 4269.6665
$ws.Range("L74").Value = 1853
$ws.Range("M74").Value = -3395.6665
$ws.Range("N74").Value = -3601
$ws.Range("H77").Value = 4049.9697
$ws.Range("I77").Value = 4269.6665
$ws.Range("J77").Value = 1853
$ws.Range("K77").Value = 21348.3325
$ws.Range("L77").Value = 9265
$ws.Range("M77").Value = -16980.3325
$ws.Range("N77").Value = -18001
$ws.Range("H119").Value = 29000
$ws.Range("J119").Value = 29000
$ws.Range("L119").Value = 29000
$ws.Range("N119").Value = -38676
$ws.Range("H136").Value = 2476.9524
$ws.Range("I136").Value = 1393
$ws.Range("J136").Value = 4644.857
$ws.Range("K136").Value = 4179
$ws.Range("L136").Value = 13934.571
$ws.Range("M136").Value = -1629
$ws.Range("N136").Value = -19034.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2933.2222
$ws.Range("I107").Value = 3055.1428
$ws.Range("J107").Value = 2506.5
$ws.Range("K107").Value = 3055.1428
$ws.Range("L107").Value = 2506.5
$ws.Range("M107").Value = -1135.1428
$ws.Range("N107").Value = -6346.5
$ws.Range("H134").Value = 3770.3684
$ws.Range("I134").Value = 3753.6155
$ws.Range("J134").Value = 3806.6667
$ws.Range("K134").Value = 11260.8465
$ws.Range("L134").Value = 11420.0001
$ws.Range("M134").Value = -8725.8465
$ws.Range("N134").Value = -16490.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1991.8334
$ws.Range("I31").Value = 1465.2322
$ws.Range("K31").Value = 1465.2322
$ws.Range("M31").Value = -1170.2322
$ws.Range("H34").Value = 1991.8334
$ws.Range("I34").Value = 1465.2322
$ws.Range("K34").Value = 1465.2322
$ws.Range("M34").Value = -1263.2322
$ws.Range("H107").Value = 607.625
$ws.Range("I107").Value = 417.26666
$ws.Range("J107").Value = 924.8889
$ws.Range("K107").Value = 417.26666
$ws.Range("L107").Value = 924.8889
$ws.Range("M107").Value = 1502.73334
$ws.Range("N107").Value = -4764.8889
$ws.Range("H132").Value = 2973.4167
$ws.Range("I132").Value = 1074.6666
$ws.Range("J132").Value = 3606.3333
$ws.Range("K132").Value = 3223.9998
$ws.Range("L132").Value = 10818.9999
$ws.Range("M132").Value = -693.9998000000001
$ws.Range("N132").Value = -15878.9999
$ws.Range("H134").Value = 3706.7778
$ws.Range("I134").Value = 1670.125
$ws.Range("K134").Value = 5010.375
$ws.Range("M134").Value = -2475.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 43.714287
$ws.Range("I14").Value = 43.714287
$ws.Range("K14").Value = 131.142861
$ws.Range("M14").Value = 41.85713900000002
$ws.Range("H33").Value = 85
$ws.Range("I33").Value = 83.333336
$ws.Range("J33").Value = 90
$ws.Range("K33").Value = 500.000016
$ws.Range("L33").Value = 540
$ws.Range("M33").Value = -217.000016
$ws.Range("N33").Value = -1106
$ws.Range("H121").Value = 11111991
$ws.Range("I121").Value = 612.25
$ws.Range("J121").Value = 20001094
$ws.Range("K121").Value = 1836.75
$ws.Range("L121").Value = 60003282
$ws.Range("M121").Value = -526.75
$ws.Range("N121").Value = -60005902
$ws.Range("H129").Value = 5137.7036
$ws.Range("I129").Value = 1920.8
$ws.Range("J129").Value = 7030
$ws.Range("K129").Value = 5762.4
$ws.Range("L129").Value = 21090
$ws.Range("M129").Value = -762.3999999999996
$ws.Range("N129").Value = -31090
$ws.Range("H131").Value = 2192.111
$ws.Range("I131").Value = 459.41177
$ws.Range("J131").Value = 2727.6726
$ws.Range("K131").Value = 1378.23531
$ws.Range("L131").Value = 8183.0178
$ws.Range("M131").Value = 3661.76469
$ws.Range("N131").Value = -18263.0178

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 18875
$ws.Range("J121").Value = 18875
$ws.Range("L121").Value = 18875
$ws.Range("N121").Value = -22369
$ws.Range("H132").Value = 3243.4736
$ws.Range("I132").Value = 2892.4
$ws.Range("J132").Value = 3633.5557
$ws.Range("K132").Value = 8677.200000000001
$ws.Range("L132").Value = 10900.6671
$ws.Range("M132").Value = -6147.200000000001
$ws.Range("N132").Value = -15960.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 77693020
$ws.Range("I22").Value = 126250440
$ws.Range("J22").Value = 1158
$ws.Range("K22").Value = 126250440
$ws.Range("L22").Value = 1158
$ws.Range("M22").Value = -126250145
$ws.Range("N22").Value = -1748
$ws.Range("H27").Value = 77693020
$ws.Range("I27").Value = 126250440
$ws.Range("J27").Value = 1158
$ws.Range("K27").Value = 126250440
$ws.Range("L27").Value = 1158
$ws.Range("M27").Value = -126250333
$ws.Range("N27").Value = -1372
$ws.Range("H100").Value = 250005000
$ws.Range("I100").Value = 10000
$ws.Range("K100").Value = 10000
$ws.Range("M100").Value = -9459
$ws.Range("H119").Value = 25834.285
$ws.Range("J119").Value = 25834.285
$ws.Range("L119").Value = 25834.285
$ws.Range("N119").Value = -35510.285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 10417147
$ws.Range("I107").Value = 12820859
$ws.Range("K107").Value = 38462577
$ws.Range("M107").Value = -38460657
$ws.Range("H119").Value = 29079.2
$ws.Range("J119").Value = 29079.2
$ws.Range("L119").Value = 29079.2
$ws.Range("N119").Value = -38755.2
$ws.Range("H132").Value = 2821.125
$ws.Range("I132").Value = 1307.7858
$ws.Range("K132").Value = 3923.3574
$ws.Range("M132").Value = -1393.3574
$ws.Range("H136").Value = 1784.4736
$ws.Range("I136").Value = 1429.6774
$ws.Range("J136").Value = 3355.7144
$ws.Range("K136").Value = 4289.0322
$ws.Range("L136").Value = 10067.1432
$ws.Range("M136").Value = -1739.0322
$ws.Range("N136").Value = -15167.1432
